$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H17").Value = 2499
$ws.Range("J17").Value = 2499
$ws.Range("L17").Value = 7497
$ws.Range("N17").Value = -7833
$ws.Range("H18").Value = 2113.2
$ws.Range("H80").Value = 1387.8
$ws.Range("I80").Value = 1196.6666
$ws.Range("K80").Value = 3589.9998
$ws.Range("M80").Value = -2591.9998
$ws.Range("H83").Value = 1387.8
$ws.Range("I83").Value = 1196.6666
$ws.Range("K83").Value = 10769.9994
$ws.Range("M83").Value = -5777.999400000001
$ws.Range("H86").Value = 1125.1
$ws.Range("I86").Value = 801
$ws.Range("J86").Value = 1449.2
$ws.Range("K86").Value = 801
$ws.Range("L86").Value = 1449.2
$ws.Range("M86").Value = 322
$ws.Range("N86").Value = -3695.2
$ws.Range("H89").Value = 1125.1
$ws.Range("I89").Value = 801
$ws.Range("J89").Value = 1449.2
$ws.Range("K89").Value = 4005
$ws.Range("L89").Value = 7246
$ws.Range("M89").Value = 1611
$ws.Range("N89").Value = -18478
$ws.Range("H111").Value = 1112.5
$ws.Range("I111").Value = 1112.5
$ws.Range("K111").Value = 3337.5
$ws.Range("M111").Value = -270.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2557.3333
$ws.Range("I74").Value = 1470.3334
$ws.Range("K74").Value = 1470.3334
$ws.Range("M74").Value = -596.3334
$ws.Range("H77").Value = 2557.3333
$ws.Range("I77").Value = 1470.3334
$ws.Range("K77").Value = 7351.666999999999
$ws.Range("M77").Value = -2983.666999999999
$ws.Range("H95").Value = 21952.5
$ws.Range("J95").Value = 21952.5
$ws.Range("L95").Value = 21952.5
$ws.Range("N95").Value = -27444.5
$ws.Range("H102").Value = 2010
$ws.Range("I102").Value = 2010
$ws.Range("K102").Value = 2010
$ws.Range("M102").Value = -388
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H122").Value = 1716.8
$ws.Range("I122").Value = 1716.8
$ws.Range("K122").Value = 5150.4
$ws.Range("M122").Value = -2700.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6950
$ws.Range("I20").Value = 6950
$ws.Range("K20").Value = 6950
$ws.Range("M20").Value = -6703
$ws.Range("H107").Value = 1200
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H31").Value = 2851.6428
$ws.Range("I31").Value = 1444.2858
$ws.Range("J31").Value = 4259
$ws.Range("K31").Value = 1444.2858
$ws.Range("L31").Value = 4259
$ws.Range("M31").Value = -1149.2858
$ws.Range("N31").Value = -4849
$ws.Range("H34").Value = 2851.6428
$ws.Range("I34").Value = 1444.2858
$ws.Range("J34").Value = 4259
$ws.Range("K34").Value = 1444.2858
$ws.Range("L34").Value = 4259
$ws.Range("M34").Value = -1242.2858
$ws.Range("N34").Value = -4663
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H60").Value = 420
$ws.Range("I60").Value = 318.875
$ws.Range("J60").Value = 824.5
$ws.Range("K60").Value = 956.625
$ws.Range("L60").Value = 2473.5
$ws.Range("M60").Value = -705.625
$ws.Range("N60").Value = -2975.5
$ws.Range("H64").Value = 500
$ws.Range("I64").Value = 500
$ws.Range("K64").Value = 1500
$ws.Range("M64").Value = -1230
$ws.Range("H67").Value = 500
$ws.Range("I67").Value = 500
$ws.Range("K67").Value = 1500
$ws.Range("M67").Value = -564
$ws.Range("H131").Value = 1551.091
$ws.Range("I131").Value = 1177
$ws.Range("K131").Value = 3531
$ws.Range("M131").Value = 1509
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6008.5
$ws.Range("I70").Value = 6008
$ws.Range("J70").Value = 6009
$ws.Range("K70").Value = 6008
$ws.Range("L70").Value = 6009
$ws.Range("M70").Value = -5738
$ws.Range("N70").Value = -6549
$ws.Range("H73").Value = 6008.5
$ws.Range("I73").Value = 6008
$ws.Range("J73").Value = 6009
$ws.Range("K73").Value = 6008
$ws.Range("L73").Value = 6009
$ws.Range("M73").Value = -5072
$ws.Range("N73").Value = -7881
$ws.Range("H97").Value = 904.875
$ws.Range("I97").Value = 819.8570999999999
$ws.Range("K97").Value = 819.8570999999999
$ws.Range("M97").Value = -323.8570999999999
$ws.Range("H107").Value = 641.125
$ws.Range("I107").Value = 583.6667
$ws.Range("J107").Value = 675.6
$ws.Range("K107").Value = 583.6667
$ws.Range("L107").Value = 675.6
$ws.Range("M107").Value = 1336.3333
$ws.Range("N107").Value = -4515.6
$ws.Range("H113").Value = 943.75
$ws.Range("I113").Value = 937.5
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 937.5
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 1232.5
$ws.Range("N113").Value = -5290
$ws.Range("H132").Value = 828.25
$ws.Range("I132").Value = 828.25
$ws.Range("K132").Value = 2484.75
$ws.Range("M132").Value = 45.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2799.6667
$ws.Range("I7").Value = 2849.5
$ws.Range("K7").Value = 2849.5
$ws.Range("M7").Value = -2737.5
$ws.Range("H46").Value = 482.4
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H126").Value = 2799.6667
$ws.Range("I126").Value = 2849.5
$ws.Range("K126").Value = 8548.5
$ws.Range("M126").Value = -6078.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39999.5
$ws.Range("J70").Value = 39999.5
$ws.Range("L70").Value = 39999.5
$ws.Range("N70").Value = -40629.5
$ws.Range("H73").Value = 39999.5
$ws.Range("J73").Value = 39999.5
$ws.Range("L73").Value = 39999.5
$ws.Range("N73").Value = -42183.5
